$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab22")

# --- Update refreshed statistics for "Africa, Fragile States" (row 97) ---
$ws.Range("C97").Value = 31.624337064585202
$ws.Range("D97").Value = 1.4714017087003199
$ws.Range("E97").Value = 3.44642195504687
$ws.Range("F97").Value = 0.1032730827561
$ws.Range("G97").Value = 250446.62730692601
$ws.Range("H97").Value = 14298.353025403499
$ws.Range("I97").Value = 27665.173888703699
$ws.Range("J97").Value = 962.51203225982204

# --- Update refreshed statistics for "ROW, Fragile States" (row 98) ---
$ws.Range("C98").Value = 15.561673464932399
$ws.Range("D98").Value = 1.0423373870709201
$ws.Range("E98").Value = 1.17065038347686
$ws.Range("F98").Value = 0.027339055875
$ws.Range("G98").Value = 137421.45062971799
$ws.Range("H98").Value = 9403.3788966062893
$ws.Range("I98").Value = 11013.8110996097
$ws.Range("J98").Value = 367.29606944624101

# --- Fix mis-encoded accented characters in the regional-groupings footnote ---
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = ""Community of Sahel-Saharan States"";COMESA = ""Common Market for Eastern and Southern Africa"";EAC = ""East African Community"";ECCAS = ""Economic Community of Central African States"";ECOWAS = ""Economic Community of West African States"";IGAD = ""Intergovernmental Authority on Development"";SADC = ""Southern African Development Community"";UMA = ""Arab Maghreb Union"";PALOP = ""Países Africanos de Língua Oficial Portuguesa"";ASEAN = ""Association of Southeast Asian Nations"";MERCOSUR = ""Mercado Común del Sur"".EU27 = ""European Union (27 members)"".OECD = ""Organisation for Economic Co-operation and Development"".";
